$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row above the current row 181. This pushes the
# existing rows 181-221 down to 182-222 (dimension grows from R221 to R222).
$ws.Rows.Item(181).Insert()

# Populate the newly inserted row 181 with the new weekly record.
$ws.Range("A181").Value = 10
$ws.Range("B181").Value = "Vega Modelo de Temuco"
$ws.Range("C181").Value = "La Araucanía"
$ws.Range("D181").Value = 44508
$ws.Range("E181").Value = 9
$ws.Range("F181").Value = 100112044
$ws.Range("G181").Value = "Perejil"
$ws.Range("H181").Value = "Sin especificar"
$ws.Range("I181").Value = "Primera"
$ws.Range("J181").Value = 40
$ws.Range("K181").Value = 4000
$ws.Range("L181").Value = 5000
$ws.Range("M181").Value = 4500
$ws.Range("N181").Value = "$/docena de atados (3 kilos)"
$ws.Range("O181").Value = "Provincia de Cautín"
$ws.Range("P181").Value = 1500
$ws.Range("Q181").Value = 3
$ws.Range("R181").Value = "Hortaliza"
